# This script applies the "Updated symbol list" data refresh to the cryptos worksheet:
#  - Row 48/49 Coin+Link values are swapped (BOLO <-> CoinbaseStockToken)
#  - Price (col D) and Volume(1h) (col E) values are refreshed for the affected rows
#
# All Price/Volume cells are stored as literal text in the workbook (e.g. "330.80", "1.53%")
# rather than as numbers, so every new value below is also written as text. Writing a
# numeric-looking string straight into a General-formatted cell would make Excel silently
# reinterpret it as a number/percentage, so each target cell is first switched to the Text
# number format ("@") and only then assigned its new string value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin / Link text swap (rows 48-49) ---------------------------------------------------
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"

# --- Price / Volume(1h) refresh ------------------------------------------------------------
$priceVolumeUpdates = [ordered]@{
    "D2" = "330.39"
    "E2" = "1.20%"
    "D3" = "44.37"
    "E3" = "0.47%"
    "D4" = "5.492"
    "E4" = "-1.20%"
    "D5" = "0.08027"
    "E5" = "0.09%"
    "D6" = "2.098"
    "E6" = "11.27%"
    "D7" = "2.657"
    "E7" = "2.18%"
    "D8" = "0.9534"
    "E8" = "1.11%"
    "D9" = "0.1145"
    "E9" = "-1.11%"
    "D10" = "0.1906"
    "E10" = "4.21%"
    "D11" = "10.20"
    "E11" = "18.42%"
    "D12" = "0.1002"
    "E12" = "2.96%"
    "D13" = "0.04802"
    "E13" = "9.99%"
    "E14" = "0.09%"
    "D15" = "0.001275"
    "E15" = "0.34%"
    "D16" = "0.04080"
    "E16" = "-3.29%"
    "D17" = "0.005890"
    "E17" = "-1.82%"
    "D18" = "3.368"
    "E18" = "-6.50%"
    "D19" = "4.399"
    "E19" = "2.39%"
    "E20" = "-1.09%"
    "D21" = "0.1382"
    "E21" = "0.16%"
    "E22" = "-2.75%"
    "E23" = "1.80%"
    "D24" = "0.004355"
    "E24" = "-2.99%"
    "D25" = "0.0001201"
    "E25" = "-4.87%"
    "D26" = "0.0003742"
    "E26" = "-6.34%"
    "D38" = "0.02592"
    "E38" = "-0.63%"
    "D39" = "0.05789"
    "E39" = "7.30%"
    "D40" = "0.007584"
    "E40" = "-0.16%"
    "D41" = "0.1407"
    "E41" = "1.04%"
    "D42" = "0.007157"
    "E42" = "-1.69%"
    "E43" = "-0.23%"
    "D44" = "0.009064"
    "E44" = "2.72%"
    "D45" = "0.00006999"
    "E45" = "1.03%"
    "E46" = "-0.12%"
    "D47" = "0.0005798"
    "E47" = "-0.23%"
    "D48" = "0.003501"
    "E48" = "-1.69%"
    "D49" = "0.003529"
    "E49" = "55.18%"
    "D50" = "0.00002101"
    "E50" = "-0.12%"
    "D51" = "0.0002001"
    "E51" = "-0.12%"
}

foreach ($cellRef in $priceVolumeUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceVolumeUpdates[$cellRef]
}
